# Refresh the "cryptos" price/volume snapshot (GitHub Actions scrape update).
#
# Each row in the table below carries the new "Price" (column D) and/or
# "Volume(1h)" (column E) text for that sheet row. Both columns store plain
# text in this workbook (e.g. "62.857.00", "  +1.54%  "), so for column D we
# briefly force a Text number format before writing the value and then clear
# the formatting again - this stops Excel's automatic type inference from
# reinterpreting decimal-look-alike strings (like "1.00" or "567.23") as
# numbers, while leaving the cell style exactly as it was before the edit.
# Column E's values keep their surrounding spaces and already round-trip as
# text untouched, so no such trick is needed there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 2; D = "62.857.00"; E = "  +1.54%  " },
    @{ Row = 3; D = "2.444.63"; E = "  +1.88%  " },
    @{ Row = 5; D = "567.23"; E = "  +1.00%  " },
    @{ Row = 6; D = "146.04"; E = "  +2.65%  " },
    @{ Row = 7; D = "1.00" },
    @{ Row = 8; E = "  +0.30%  " },
    @{ Row = 9; E = "  +2.88%  " },
    @{ Row = 10; E = "  +0.54%  " },
    @{ Row = 11; D = "5.30"; E = "  +0.91%  " },
    @{ Row = 12; E = "  +1.98%  " },
    @{ Row = 13; D = "27.04"; E = "  +5.98%  " },
    @{ Row = 14; E = "  +6.40%  " },
    @{ Row = 15; D = "2.799.55"; E = "  -1.17%  " },
    @{ Row = 16; D = "62.569.50"; E = "  +1.21%  " },
    @{ Row = 17; D = "2.432.38"; E = "  +1.90%  " },
    @{ Row = 18; D = "11.29"; E = "  +0.79%  " },
    @{ Row = 19; D = "6.95"; E = "  +2.52%  " },
    @{ Row = 20; D = "324.50"; E = "  +1.21%  " },
    @{ Row = 21; E = "  +1.23%  " },
    @{ Row = 22; D = "1.00"; E = "  +0.09%  " },
    @{ Row = 23; D = "1.86"; E = "  +7.06%  " },
    @{ Row = 24; D = "67.37"; E = "  +2.19%  " },
    @{ Row = 25; E = "  -0.79%  " },
    @{ Row = 26; D = "592.61"; E = "  +5.77%  " },
    @{ Row = 27; E = "  +9.65%  " },
    @{ Row = 28; D = "2.560.69"; E = "  +1.61%  " },
    @{ Row = 29; D = "8.47"; E = "  +4.14%  " },
    @{ Row = 30; D = "0.999"; E = "  -0.18%  " },
    @{ Row = 31; E = "  +5.44%  " },
    @{ Row = 32; E = "  +0.88%  " },
    @{ Row = 33; E = "  +0.61%  " },
    @{ Row = 34; E = "  +3.45%  " },
    @{ Row = 35; D = "4.88"; E = "  +4.98%  " },
    @{ Row = 36; E = "  -0.15%  " },
    @{ Row = 37; E = "  +1.71%  " },
    @{ Row = 38; D = "5.45"; E = "  +0.97%  " },
    @{ Row = 39; D = "18.83"; E = "  +1.62%  " },
    @{ Row = 40; E = "  -2.53%  " },
    @{ Row = 41; E = "  +2.69%  " },
    @{ Row = 42; E = "  +0.47%  " },
    @{ Row = 43; D = "2.45"; E = "  +9.23%  " },
    @{ Row = 44; D = "149.30"; E = "  +1.54%  " },
    @{ Row = 45; E = "  +2.43%  " },
    @{ Row = 47; D = "20.62"; E = "  +4.58%  " },
    @{ Row = 48; E = "  +3.05%  " },
    @{ Row = 50; E = "  +1.11%  " },
    @{ Row = 51; E = "  +3.97%  " }
)

foreach ($chg in $changes) {
    $row = $chg.Row

    if ($chg.ContainsKey("D")) {
        $priceCell = $ws.Cells.Item($row, 4)
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $chg.D
        $priceCell.ClearFormats()
    }

    if ($chg.ContainsKey("E")) {
        $volumeCell = $ws.Cells.Item($row, 5)
        $volumeCell.Value = $chg.E
    }
}
